# Add a new "2022-Q1" sheet right before the "总计" (Total) sheet, carrying
# the per-fund holding detail, and insert a corresponding summary row at the
# top of the "总计" sheet's data.

$wb = $excel.ActiveWorkbook

# A cell that already carries the workbook's "header / index column" style
# (bold font + border + centered) - used as a formatting template below so
# the new cells match the look of the existing quarterly sheets.
$styleTemplate = $wb.Worksheets.Item("2021-Q4").Range("B1")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Header row (row 1) - matches the other quarterly detail sheets.
$styleTemplate.Copy($newSheet.Range("B1"))
$newSheet.Range("B1").Value = "基金代码"
$styleTemplate.Copy($newSheet.Range("C1"))
$newSheet.Range("C1").Value = "基金名称"
$styleTemplate.Copy($newSheet.Range("D1"))
$newSheet.Range("D1").Value = "基金规模"
$styleTemplate.Copy($newSheet.Range("E1"))
$newSheet.Range("E1").Value = "股票总仓位"
$styleTemplate.Copy($newSheet.Range("F1"))
$newSheet.Range("F1").Value = "仓位占比"
$styleTemplate.Copy($newSheet.Range("G1"))
$newSheet.Range("G1").Value = "持有市值(亿元)"
$styleTemplate.Copy($newSheet.Range("H1"))
$newSheet.Range("H1").Value = "仓位排名"

# Data row 2.
$styleTemplate.Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "513030"
$newSheet.Range("B2").Style = "Normal"

$newSheet.Range("C2").Value = "华安国际龙头(DAX)ETFQDII"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "6.49"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "92.80"
$newSheet.Range("E2").Style = "Normal"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "6.92"
$newSheet.Range("F2").Style = "Normal"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.4491"
$newSheet.Range("G2").Style = "Normal"

$newSheet.Range("H2").Value = 3

# Data row 3.
$styleTemplate.Copy($newSheet.Range("A3"))
$newSheet.Range("A3").Value = 1

$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "006282"
$newSheet.Range("B3").Style = "Normal"

$newSheet.Range("C3").Value = "上投摩根欧洲动力策略股票（QDII）"

$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.48"
$newSheet.Range("D3").Style = "Normal"

$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "89.68"
$newSheet.Range("E3").Style = "Normal"

$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "1.95"
$newSheet.Range("F3").Style = "Normal"

$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0094"
$newSheet.Range("G3").Style = "Normal"

$newSheet.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. Insert a new top data row in "总计" for 2022-Q1, pushing the rest down.
#    (re-fetch the sheet reference by name - the handle captured before the
#    Add() above can point at the wrong sheet once the tab order shifts.)
#    Rows are shifted by copying bottom-to-top so each existing row keeps
#    its original formatting (a plain Range.Insert blends in formatting
#    from whichever row it borrows from).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A6:D6").Copy($totalSheet.Range("A7:D7"))
$totalSheet.Range("A5:D5").Copy($totalSheet.Range("A6:D6"))
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.46

# Renumber the helper index column (A) for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
